$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) becomes bold ---
$ws.Range("A1:E1").Font.Bold = $true

# --- Fix typo in translation for "abgewöhnen" (row 12): "se déhabituer" -> "se déshabituer" ---
$ws.Range("B12").Value = "se déshabituer"

# --- New translation added for "abfinden" (row 7): "se résigner" in column C ---
$ws.Range("C7").Value = "se résigner"

# --- Scroll/selection state: select header row, scroll so row 10 is at top ---
$ws.Range("A1:E1").Select()
$excel.ActiveWindow.ScrollRow = 10
